$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A171").Value = "test"
